# Auto-generated Excel COM-interop script to apply crypto price/volume update
# (commit: "Updated cryptos list on Wed Dec 13 10:30:18 UTC 2023 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "41.166.55"
$ws.Range("E2").Value = "  -1.43%  "

# Row 3
$ws.Range("D3").Value = "2.174.13"
$ws.Range("E3").Value = "  -2.23%  "

# Row 4
$ws.Range("E4").Value = "  -0.06%  "

# Row 5
$ws.Range("D5").Value = "'250.06"
$ws.Range("E5").Value = "  -0.34%  "

# Row 6
$ws.Range("D6").Value = "'0.612"
$ws.Range("E6").Value = "  -2.84%  "

# Row 7
$ws.Range("D7").Value = "'66.10"
$ws.Range("E7").Value = "  -8.04%  "

# Row 8
$ws.Range("E8").Value = "  +0.06%  "

# Row 9
$ws.Range("D9").Value = "'0.577"
$ws.Range("E9").Value = "  -3.23%  "

# Row 10
$ws.Range("D10").Value = "'58.90"
$ws.Range("E10").Value = "  +1.12%  "

# Row 11
$ws.Range("D11").Value = "'36.21"
$ws.Range("E11").Value = "  -12.35%  "

# Row 12
$ws.Range("D12").Value = "'0.0933"
$ws.Range("E12").Value = "  -3.05%  "

# Row 13
$ws.Range("E13").Value = "  -1.67%  "

# Row 14
$ws.Range("D14").Value = "'6.82"
$ws.Range("E14").Value = "  -4.84%  "

# Row 15
$ws.Range("D15").Value = "2.504.23"
$ws.Range("E15").Value = "  -1.97%  "

# Row 16
$ws.Range("D16").Value = "'14.26"
$ws.Range("E16").Value = "  -4.38%  "

# Row 17
$ws.Range("D17").Value = "'0.845"
$ws.Range("E17").Value = "  -2.36%  "

# Row 18
$ws.Range("D18").Value = "2.155.63"
$ws.Range("E18").Value = "  -2.91%  "

# Row 19
$ws.Range("D19").Value = "41.105.05"
$ws.Range("E19").Value = "  -1.60%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0946"
$ws.Range("E20").Value = "  -2.16%  "

# Row 21
$ws.Range("D21").Value = "'71.58"
$ws.Range("E21").Value = "  -1.91%  "

# Row 22
$ws.Range("D22").Value = "'6.05"
$ws.Range("E22").Value = "  -2.81%  "

# Row 23
$ws.Range("D23").Value = "'229.93"
$ws.Range("E23").Value = "  -2.24%  "

# Row 24
$ws.Range("D24").Value = "'2.02"
$ws.Range("E24").Value = "  -5.04%  "

# Row 25
$ws.Range("D25").Value = "'3.79"
$ws.Range("E25").Value = "  -6.46%  "

# Row 26
$ws.Range("D26").Value = "'11.43"
$ws.Range("E26").Value = "  +6.54%  "

# Row 27
$ws.Range("E27").Value = "  +0.14%  "

# Row 28
$ws.Range("D28").Value = "'2.41"
$ws.Range("E28").Value = "  -5.10%  "

# Row 29
$ws.Range("D29").Value = "'168.37"

# Row 30
$ws.Range("D30").Value = "'2.02"
$ws.Range("E30").Value = "  -4.19%  "

# Row 31
$ws.Range("D31").Value = "'20.16"
$ws.Range("E31").Value = "  -3.01%  "

# Row 32
$ws.Range("D32").Value = "'0.122"
$ws.Range("E32").Value = "  -3.12%  "

# Row 33
$ws.Range("D33").Value = "'5.75"
$ws.Range("E33").Value = "  +2.24%  "

# Row 34
$ws.Range("D34").Value = "'0.0748"
$ws.Range("E34").Value = "  +2.09%  "

# Row 35
$ws.Range("E35").Value = "  -3.33%  "

# Row 36
$ws.Range("D36").Value = "'4.50"
$ws.Range("E36").Value = "  -4.94%  "

# Row 37
$ws.Range("D37").Value = "'3.93"
$ws.Range("E37").Value = "  -1.78%  "

# Row 38
$ws.Range("D38").Value = "'24.30"
$ws.Range("E38").Value = "  -5.80%  "

# Row 39
$ws.Range("D39").Value = "'0.0304"
$ws.Range("E39").Value = "  +0.69%  "

# Row 40
$ws.Range("D40").Value = "'2.21"
$ws.Range("E40").Value = "  -3.58%  "

# Row 41
$ws.Range("D41").Value = "'5.37"
$ws.Range("E41").Value = "  +9.42%  "

# Row 42
$ws.Range("D42").Value = "'5.47"
$ws.Range("E42").Value = "  -8.12%  "

# Row 43
$ws.Range("D43").Value = "'11.31"
$ws.Range("E43").Value = "  -6.74%  "

# Row 44
$ws.Range("D44").Value = "'60.50"
$ws.Range("E44").Value = "  -9.60%  "

# Row 45
$ws.Range("D45").Value = "'8.49"
$ws.Range("E45").Value = "  -3.22%  "

# Row 46
$ws.Range("D46").Value = "'0.0993"
$ws.Range("E46").Value = "  -2.93%  "

# Row 47
$ws.Range("B47").Value = "BinanceUSD"
$ws.Range("C47").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D47").Value = "'1.00"
$ws.Range("E47").Value = "  +0.03%  "

# Row 48
$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D48").Value = "'0.188"
$ws.Range("E48").Value = "  -7.03%  "

# Row 49
$ws.Range("E49").Value = "  -2.39%  "

# Row 50
$ws.Range("E50").Value = "  -4.80%  "

# Row 51
$ws.Range("E51").Value = "  -10.84%  "
